$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# Helper: write a text-typed value (e.g. "0" or "***.*") into a cell while
# keeping it on the same "label" style (s=13) used throughout this sheet for
# suppressed/undefined figures. We first force the literal string in via the
# apostrophe-prefix trick (so Excel doesn't silently re-parse "0" back into a
# number), then paste-special the number format from a known-good donor cell
# that already carries style 13, which also clears the quote-prefix styling
# that the literal-text entry would otherwise introduce.
function Set-LabelText {
    param($targetAddr, $donorAddr, $text)
    $ws.Range($targetAddr).Value = "'" + $text
    $ws.Range($donorAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial($xlPasteFormats)
}

# Donor cells (row 14 is untouched by this edit and already has style 13):
#   C14 -> text "0"
#   E14 -> text "***.*"
$donorZero = "C14"
$donorStar = "E14"

# ---------------------------------------------------------------------
# Shared strings: volume number and report week dates (rich-text runs)
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "12"
$ws.Range("C9").Characters(27, 9).Text = "3/17/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/23/2025"

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-LabelText "C16" $donorZero "0"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -83.333333333333
$ws.Range("M16").Value = -60
$ws.Range("N16").Value = -91.666666666666

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-LabelText "C17" $donorZero "0"

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 5
$ws.Range("G19").Value = 7
$ws.Range("H19").Value = -28.571428571428
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 9
$ws.Range("L19").Value = 125
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = -18.181818181818

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("F21").Value = 8
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = -38.461538461538
$ws.Range("I21").Value = 15
$ws.Range("J21").Value = 21
$ws.Range("K21").Value = -28.571428571428
$ws.Range("L21").Value = 66.666666666666
$ws.Range("M21").Value = 36.363636363636
$ws.Range("N21").Value = -68.75

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-LabelText "D24" $donorZero "0"
Set-LabelText "E24" $donorStar "***.*"
$ws.Range("L24").Value = -88.888888888888

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
Set-LabelText "D26" $donorZero "0"
Set-LabelText "E26" $donorStar "***.*"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("L26").Value = -40
$ws.Range("M26").Value = -57.142857142857

$excel.CutCopyMode = $false
